# Updates the "cryptos" sheet with refreshed price / 1h-volume-change figures,
# and swaps the Cronos / InjectiveProtocol rows (46-47) to reflect their new rank order.
# D-column price cells that look numeric ("0.636", "246.45", ...) are written via a
# temporary text NumberFormat ("@") so Excel keeps them as text instead of silently
# converting them to real numbers; the style is reset to "Normal" right afterwards so
# no stray cell formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.608.58"
$ws.Range("E2").Value = "  -2.92%  "
$ws.Range("D3").Value = "1.984.03"
$ws.Range("E3").Value = "  -3.63%  "
$ws.Range("E4").Value = "  +0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "246.45"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.17%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.636"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -4.66%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "58.50"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +7.24%  "
$ws.Range("E8").Value = "  -0.02%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "58.71"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.76%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.362"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.85%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0737"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("E12").Value = "  -2.66%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.957"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.25%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "14.63"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("D15").Value = "2.270.79"
$ws.Range("E15").Value = "  -3.75%  "
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("D17").Value = "1.959.94"
$ws.Range("E17").Value = "  -4.83%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "18.41"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +7.79%  "
$ws.Range("D19").Value = "35.527.09"
$ws.Range("E19").Value = "  -2.88%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "71.44"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").Value = "0.0₃0849"
$ws.Range("E21").Value = "  -1.68%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.23"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.70%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "232.78"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.07%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.58"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +20.54%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.25"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -5.29%  "
$ws.Range("E27").Value = "  +0.19%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "9.16"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -2.05%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "19.24"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -4.57%  "
$ws.Range("E30").Value = "  -2.46%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.89"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.99%  "
$ws.Range("E32").Value = "  -6.77%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0945"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +13.32%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0596"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("E35").Value = "  +9.39%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.36"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -3.54%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  -3.78%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "5.36"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +8.19%  "
$ws.Range("E40").Value = "  -2.04%  "
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("E42").Value = "  -1.36%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "7.85"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.67%  "
$ws.Range("E44").Value = "  -1.29%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "93.47"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0901"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "16.17"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("D48").Value = "1.378.40"
$ws.Range("E48").Value = "  -2.37%  "
$ws.Range("E49").Value = "  -0.56%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "47.32"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.32%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.28"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.17%  "
